$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay text so values such as thousand-dot
# separated numbers, leading/trailing zeros, and scientific-looking
# small decimals are preserved exactly as plain text, matching the
# original inline-string cell contents.
$ws.Range("D2:D51").NumberFormat = "@"

# Update Price (column D) and Volume(1h) (column E) values for the refreshed crypto data
$ws.Range("D2").Value = "24.354.48"
$ws.Range("E2").Value = "  -3.81%  "
$ws.Range("D3").Value = "1.642.99"
$ws.Range("E3").Value = "  -5.91%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "0.9984"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "305.14"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("D7").Value = "0.3624"
$ws.Range("E7").Value = "  -5.20%  "
$ws.Range("D8").Value = "47.31"
$ws.Range("E8").Value = "  -4.52%  "
$ws.Range("E9").Value = "  -8.66%  "
$ws.Range("D10").Value = "1.111"
$ws.Range("E10").Value = "  -8.80%  "
$ws.Range("D11").Value = "0.06912"
$ws.Range("E11").Value = "  -9.47%  "
$ws.Range("D12").Value = "0.9982"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "5.939"
$ws.Range("E13").Value = "  -8.39%  "
$ws.Range("D14").Value = "19.08"
$ws.Range("E14").Value = "  -10.62%  "
$ws.Range("D15").Value = "1.648.34"
$ws.Range("E15").Value = "  -5.71%  "
$ws.Range("D16").Value = "6.527"
$ws.Range("E16").Value = "  -8.01%  "
$ws.Range("D17").Value = "0.00001043"
$ws.Range("E17").Value = "  -9.37%  "
$ws.Range("D18").Value = "0.06473"
$ws.Range("E18").Value = "  -4.30%  "
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "76.78"
$ws.Range("E20").Value = "  -10.78%  "
$ws.Range("D21").Value = "5.886"
$ws.Range("E21").Value = "  -9.62%  "
$ws.Range("D22").Value = "15.69"
$ws.Range("E22").Value = "  -10.72%  "
$ws.Range("D23").Value = "12.10"
$ws.Range("E23").Value = "  -7.80%  "
$ws.Range("D24").Value = "24.371.14"
$ws.Range("E24").Value = "  -3.61%  "
$ws.Range("D25").Value = "2.407"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").Value = "2.329"
$ws.Range("E26").Value = "  -18.77%  "
$ws.Range("D27").Value = "145.32"
$ws.Range("E27").Value = "  -5.58%  "
$ws.Range("E28").Value = "  -11.68%  "
$ws.Range("D29").Value = "1.828.51"
$ws.Range("E29").Value = "  -5.85%  "
$ws.Range("D30").Value = "124.56"
$ws.Range("E30").Value = "  -6.76%  "
$ws.Range("D31").Value = "1.145"
$ws.Range("E31").Value = "  -4.98%  "
$ws.Range("D32").Value = "4.047"
$ws.Range("E32").Value = "  -4.63%  "
$ws.Range("D33").Value = "5.594"
$ws.Range("E33").Value = "  -21.05%  "
# Rows 34 and 35 swap content (Stellar moves up to row 34, WEMIXTOKEN moves to row 35)
# along with updated Price/Volume values.
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "0.08308"
$ws.Range("E34").Value = "  -5.33%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "1.677"
$ws.Range("E35").Value = "  -7.42%  "
$ws.Range("D36").Value = "12.26"
$ws.Range("E36").Value = "  -13.94%  "
$ws.Range("D37").Value = "5.156"
$ws.Range("E37").Value = "  -10.75%  "
$ws.Range("D38").Value = "0.06031"
$ws.Range("E38").Value = "  -10.10%  "
$ws.Range("D39").Value = "0.02202"
$ws.Range("E39").Value = "  -11.99%  "
$ws.Range("D40").Value = "8.204"
$ws.Range("E40").Value = "  -11.93%  "
$ws.Range("D41").Value = "1.200"
$ws.Range("E41").Value = "  -6.53%  "
$ws.Range("D42").Value = "0.2031"
$ws.Range("E42").Value = "  -10.10%  "
$ws.Range("D43").Value = "0.9976"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "0.5832"
$ws.Range("E44").Value = "  -11.38%  "
$ws.Range("D45").Value = "3.712"
$ws.Range("E45").Value = "  -4.72%  "
$ws.Range("D46").Value = "12.65"
$ws.Range("E46").Value = "  -11.34%  "
$ws.Range("D47").Value = "0.5574"
$ws.Range("E47").Value = "  -11.28%  "
$ws.Range("D48").Value = "121.36"
$ws.Range("E48").Value = "  -7.72%  "
$ws.Range("D49").Value = "1.929"
$ws.Range("E49").Value = "  -11.77%  "
$ws.Range("D50").Value = "0.06870"
$ws.Range("E50").Value = "  -7.11%  "
$ws.Range("D51").Value = "73.63"
$ws.Range("E51").Value = "  -8.89%  "
